$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - renamed
$ws.Range("A1").Value = "Value1_Max"
$ws.Range("B1").Value = "Value2_Max"
$ws.Range("C1").Value = "Value3_Max"
$ws.Range("D1").Value = "Value4_Max"
$ws.Range("E1").Value = "Value5_Max"
$ws.Range("F1").Value = "Value6_Max"
$ws.Range("G1").Value = " Value7_Max"
$ws.Range("H1").Value = "Value1_Warning"
$ws.Range("I1").Value = "Value2_Warning"
$ws.Range("J1").Value = "Value3_Warning"
$ws.Range("K1").Value = "Value4_Warning"
$ws.Range("L1").Value = "Value5_Warning"
$ws.Range("M1").Value = "Value6_Warning"
$ws.Range("N1").Value = "Value7_Warning"

# Row 2: A2, B2, F2, G2, C2, D2, E2 stay as-is (already correct: 0.5 / 80 / 0.004)

# H2, I2, M2, N2 become text values (quote-prefixed so they stay text, not numbers)
$ws.Range("H2").Value = "'0.3"
$ws.Range("I2").Value = "'0.3"
$ws.Range("M2").Value = "'0.002"
$ws.Range("N2").Value = "'0.002"

# J2, K2, L2 updated numeric values 40 -> 60
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 60

# Update selection to O2
$ws.Range("O2").Select()
